$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Projet")
$ws2 = $wb.Worksheets.Item("Cahier de test")

$ws2.Range("B4").Value = 0
$ws2.Range("B5").Value = 1
$ws2.Range("B6").Value = 2
$ws2.Range("B7").Value = 0
$ws2.Range("B8").Value = 1
$ws2.Range("B9").Value = 0
$ws2.Range("B10").Value = 0
$ws2.Range("B11").Value = 1
$ws2.Range("B12").Value = 1
$ws2.Range("B13").Value = 0

$ws1.Activate()
$ws1.Range("I5").Select()
$ws2.Activate()
$ws2.Range("C5").Select()
